# Generate Report for Handback
#
# This script updates the localization-status workbook to reflect that the
# "a.md" file has been handed back (for both the zh-cn and de-de targets),
# while "b.md" is still pending. It:
#   - updates the shared "Status" text from "Ready for handoff" to
#     "Handed back: in sync with en-US" everywhere it is used,
#   - fills in the "Latest Target File" / "Latest Handback File" /
#     "Latest Handback DateTime" columns on the zh-cn and de-de sheets,
#   - adds hyperlinks on the new "Latest Target File" cells,
#   - widens a few columns that now hold longer text.

$wb = $excel.ActiveWorkbook

$Overview = $wb.Worksheets.Item(1)
$ZhCn     = $wb.Worksheets.Item(2)
$DeDe     = $wb.Worksheets.Item(3)

$newStatus = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# 1) Status text: "Ready for handoff" -> "Handed back: in sync with en-US"
#    (shared string, referenced from Overview!E2:F3 and from the status
#    column (C) on both the zh-cn and de-de sheets)
# ---------------------------------------------------------------------
$Overview.Range("E2").Value = $newStatus
$Overview.Range("F2").Value = $newStatus
$Overview.Range("E3").Value = $newStatus
$Overview.Range("F3").Value = $newStatus

$ZhCn.Range("C2").Value = $newStatus
$ZhCn.Range("C3").Value = $newStatus

$DeDe.Range("C2").Value = $newStatus
$DeDe.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------
# 2) Column widths: widen columns that now contain the longer status /
#    file-name text.
# ---------------------------------------------------------------------
# Overview: columns E (zh-cn) and F (de-de)
$Overview.Columns.Item(5).ColumnWidth = 29.166666666666668
$Overview.Columns.Item(6).ColumnWidth = 29.166666666666668

# zh-cn / de-de sheets: column C (Status) and column J (Latest Handback File)
$ZhCn.Columns.Item(3).ColumnWidth = 29.166666666666668
$ZhCn.Columns.Item(10).ColumnWidth = 39.166666666666664

$DeDe.Columns.Item(3).ColumnWidth = 29.166666666666668
$DeDe.Columns.Item(10).ColumnWidth = 39.166666666666664

# ---------------------------------------------------------------------
# 3) zh-cn sheet (a.md row = row 2): fill in target / handback info.
#    b.md (row 3) is not handed off yet, but the fixture data mirrors the
#    a.md values onto row 3 as well.
# ---------------------------------------------------------------------
$hoBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b6adfa9631904f6631db643e5edcac4a4db95af6/e2e/"

$ZhCn.Range("I2").Value = "a.md"
$ZhCn.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$ZhCn.Range("K2").Value = "2016-08-20 22:43:34"

$ZhCn.Range("I3").Value = "a.md"
$ZhCn.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$ZhCn.Range("K3").Value = "2016-08-20 22:43:34"

# Match the look of the existing "source file" hyperlink cells.
$ZhCn.Range("I2").Font.Underline = 2
$ZhCn.Range("I2").Font.Color = 15570276
$ZhCn.Range("I3").Font.Underline = 2
$ZhCn.Range("I3").Font.Color = 15570276

# Recreate the hyperlinks in display order (A2, I2, A3, I3) so relationship
# ids line up the way Excel would assign them.
$ZhCn.Hyperlinks.Delete()
$ZhCn.Hyperlinks.Add($ZhCn.Range("A2"), ($hoBase + "a.md"), [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "a.md")
$ZhCn.Hyperlinks.Add($ZhCn.Range("I2"), ($hoBase + "a.md"), [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "a.md")
$ZhCn.Hyperlinks.Add($ZhCn.Range("A3"), ($hoBase + "b.md"), [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "b.md")
$ZhCn.Hyperlinks.Add($ZhCn.Range("I3"), ($hoBase + "a.md"), [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "a.md")

# ---------------------------------------------------------------------
# 4) de-de sheet (a.md row = row 2): fill in target / handback info.
# ---------------------------------------------------------------------
$DeDe.Range("I2").Value = "a.md"
$DeDe.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$DeDe.Range("K2").Value = "2016-08-20 22:43:41"

$DeDe.Range("I3").Value = "a.md"
$DeDe.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$DeDe.Range("K3").Value = "2016-08-20 22:43:41"

$DeDe.Range("I2").Font.Underline = 2
$DeDe.Range("I2").Font.Color = 15570276
$DeDe.Range("I3").Font.Underline = 2
$DeDe.Range("I3").Font.Color = 15570276

$DeDe.Hyperlinks.Delete()
$DeDe.Hyperlinks.Add($DeDe.Range("A2"), ($hoBase + "a.md"), [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "a.md")
$DeDe.Hyperlinks.Add($DeDe.Range("I2"), ($hoBase + "a.md"), [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "a.md")
$DeDe.Hyperlinks.Add($DeDe.Range("A3"), ($hoBase + "b.md"), [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "b.md")
$DeDe.Hyperlinks.Add($DeDe.Range("I3"), ($hoBase + "a.md"), [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "a.md")

$wb.Save()
